# autoast/alan_jobs.xlsx update
#
# Commit message:
#   - moved skip if complete and failed logic one indent inside the for loop
#   - deleted autoast_indexing after copying it and naming it V2 Cuisinart_Dev.py
#
# The accompanying job-tracking workbook (ast_config sheet) reflects the
# effect of that logic change on the in-flight job in row 2:
#   - J2 (suppress_map_creation) goes from TRUE to "false"
#   - M2 (ast_condition) goes from "COMPLETE" to "Queued"
#
# Both new values are written via Copy / PasteSpecial (values only) so
# that:
#   * "false" ends up stored as plain text (matching the other text
#     cells in the same row, e.g. K2) instead of being auto-coerced into
#     a Boolean the way a direct .Value = "false" assignment would be,
#     and
#   * the destination cells keep their own existing cell formatting,
#     since PasteSpecial values-only does not touch the destination
#     style (in particular M2 keeps its "quote prefix" text style).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ast_config")

# J2: suppress_map_creation -> "false" (copy the existing text "false"
# from K2 on the same row so it is written as text, not a Boolean)
$ws.Range("K2").Copy()
$ws.Range("J2").PasteSpecial(-4163)

# M2: ast_condition -> "Queued" (use a scratch cell + paste-values so
# M2 keeps its current style/number format)
$scratch = $ws.Range("Z100")
$scratch.Value = "Queued"
$scratch.Copy()
$ws.Range("M2").PasteSpecial(-4163)
$scratch.Clear()

$excel.CutCopyMode = 0
